$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report generation date in D5
$ws.Range("D5").Value = "2026.02.22 22:00"

# Append new trade history rows (296-305)
# Row 296
$ws.Range("A296").Value = "2026.02.20 09:13:50"
$ws.Range("B296").Value = 857029424
$ws.Range("C296").Value = "BTCUSD"
$ws.Range("D296").Value = "buy"
$ws.Range("E296").Value = 0.5
$ws.Range("F296").Value = 67885.42
$ws.Range("I296").Value = "2026.02.20 12:02:32"
$ws.Range("J296").Value = 67986.69
$ws.Range("K296").Value = 0
$ws.Range("L296").Value = 0
$ws.Range("M296").Value = 50.64
$ws.Range("N296").Value = "[tp 67986.69]"

# Row 297
$ws.Range("A297").Value = "2026.02.20 11:03:46"
$ws.Range("B297").Value = 858034226
$ws.Range("C297").Value = "XAUUSD"
$ws.Range("D297").Value = "buy"
$ws.Range("E297").Value = 0.26
$ws.Range("F297").Value = 5015.08
$ws.Range("I297").Value = "2026.02.20 12:08:21"
$ws.Range("J297").Value = 5036.92
$ws.Range("K297").Value = 0
$ws.Range("L297").Value = 0
$ws.Range("M297").Value = 567.84
$ws.Range("N297").Value = "closePosition"

# Row 298
$ws.Range("A298").Value = "2026.02.20 12:12:27"
$ws.Range("B298").Value = 858453115
$ws.Range("C298").Value = "BTCUSD"
$ws.Range("D298").Value = "buy"
$ws.Range("E298").Value = 0.5
$ws.Range("F298").Value = 68139.89
$ws.Range("I298").Value = "2026.02.20 14:37:17"
$ws.Range("J298").Value = 67700.04
$ws.Range("K298").Value = 0
$ws.Range("L298").Value = 0
$ws.Range("M298").Value = -219.93
$ws.Range("N298").Value = "[sl 67700.04]"

# Row 299
$ws.Range("A299").Value = "2026.02.20 15:01:41"
$ws.Range("B299").Value = 859203714
$ws.Range("C299").Value = "SP500"
$ws.Range("D299").Value = "sell"
$ws.Range("E299").Value = 5
$ws.Range("F299").Value = 6854.85
$ws.Range("I299").Value = "2026.02.20 15:34:59"
$ws.Range("J299").Value = 6852.96
$ws.Range("K299").Value = 0
$ws.Range("L299").Value = 0
$ws.Range("M299").Value = 9.45
$ws.Range("N299").Value = "closePosition"

# Row 300
$ws.Range("A300").Value = "2026.02.20 15:02:01"
$ws.Range("B300").Value = 859205782
$ws.Range("C300").Value = "NAS100"
$ws.Range("D300").Value = "sell"
$ws.Range("E300").Value = 9.9
$ws.Range("F300").Value = 24756.27
$ws.Range("I300").Value = "2026.02.20 15:35:01"
$ws.Range("J300").Value = 24752.07
$ws.Range("K300").Value = 0
$ws.Range("L300").Value = 0
$ws.Range("M300").Value = 41.58
$ws.Range("N300").Value = "closePosition"

# Row 301
$ws.Range("A301").Value = "2026.02.20 15:12:31"
$ws.Range("B301").Value = 859258095
$ws.Range("C301").Value = "USDJPY"
$ws.Range("D301").Value = "buy"
$ws.Range("E301").Value = 1
$ws.Range("F301").Value = 155.28
$ws.Range("I301").Value = "2026.02.20 15:35:03"
$ws.Range("J301").Value = 155.278
$ws.Range("K301").Value = 0
$ws.Range("L301").Value = 0
$ws.Range("M301").Value = -1.29
$ws.Range("N301").Value = "closePosition"

# Row 302
$ws.Range("A302").Value = "2026.02.21 12:21:54"
$ws.Range("B302").Value = 862813462
$ws.Range("C302").Value = "BTCUSD"
$ws.Range("D302").Value = "buy"
$ws.Range("E302").Value = 0.1
$ws.Range("F302").Value = 68250.86
$ws.Range("I302").Value = "2026.02.21 12:22:10"
$ws.Range("J302").Value = 68251.36
$ws.Range("K302").Value = 0
$ws.Range("L302").Value = 0
$ws.Range("M302").Value = 0.05
$ws.Range("N302").Value = "[tp 68251.36]"

# Row 303
$ws.Range("A303").Value = "2026.02.21 12:32:00"
$ws.Range("B303").Value = 862818851
$ws.Range("C303").Value = "BTCUSD"
$ws.Range("D303").Value = "buy"
$ws.Range("E303").Value = 0.1
$ws.Range("F303").Value = 68238.96
$ws.Range("I303").Value = "2026.02.21 12:40:05"
$ws.Range("J303").Value = 68239.46
$ws.Range("K303").Value = 0
$ws.Range("L303").Value = 0
$ws.Range("M303").Value = 0.05
$ws.Range("N303").Value = "[tp 68239.46]"

# Row 304
$ws.Range("A304").Value = "2026.02.21 16:36:07"
$ws.Range("B304").Value = 862936210
$ws.Range("C304").Value = "SOLUSD"
$ws.Range("D304").Value = "buy"
$ws.Range("E304").Value = 5
$ws.Range("F304").Value = 85.63
$ws.Range("I304").Value = "2026.02.21 17:21:50"
$ws.Range("J304").Value = 86.19
$ws.Range("K304").Value = 0
$ws.Range("L304").Value = 0
$ws.Range("M304").Value = 28
$ws.Range("N304").Value = "[tp 86.19]"

# Row 305
$ws.Range("A305").Value = "2026.02.21 14:31:06"
$ws.Range("B305").Value = 862866373
$ws.Range("C305").Value = "BNBUSD"
$ws.Range("D305").Value = "buy"
$ws.Range("E305").Value = 1
$ws.Range("F305").Value = 631.68
$ws.Range("I305").Value = "2026.02.21 19:06:15"
$ws.Range("J305").Value = 623.33
$ws.Range("K305").Value = 0
$ws.Range("L305").Value = 0
$ws.Range("M305").Value = -83.5
$ws.Range("N305").Value = "[sl 623.33]"
